$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy format of row 269 (last existing data row) for column A style (date style s="2")
$ws.Range("A269").Copy()

$ws.Range("A270").PasteSpecial(-4122)
$ws.Range("A270").Value = 44344
$ws.Range("B270").Value = 0
$ws.Range("C270").Value = 3
$ws.Range("D270").Value = 48.5201358563804

$ws.Range("A271").PasteSpecial(-4122)
$ws.Range("A271").Value = 44345
$ws.Range("B271").Value = 1
$ws.Range("C271").Value = 4
$ws.Range("D271").Value = 64.69351447517387

$ws.Range("A272").PasteSpecial(-4122)
$ws.Range("A272").Value = 44346
$ws.Range("B272").Value = 0
$ws.Range("C272").Value = 4
$ws.Range("D272").Value = 64.69351447517387

$ws.Range("A273").PasteSpecial(-4122)
$ws.Range("A273").Value = 44347
$ws.Range("B273").Value = 0
$ws.Range("C273").Value = 1
$ws.Range("D273").Value = 16.17337861879347

$ws.Range("A274").PasteSpecial(-4122)
$ws.Range("A274").Value = 44348
$ws.Range("B274").Value = 0
$ws.Range("C274").Value = 1
$ws.Range("D274").Value = 16.17337861879347

$ws.Range("A275").PasteSpecial(-4122)
$ws.Range("A275").Value = 44349
$ws.Range("B275").Value = 0
$ws.Range("C275").Value = 1
$ws.Range("D275").Value = 16.17337861879347

$ws.Range("A276").PasteSpecial(-4122)
$ws.Range("A276").Value = 44350
$ws.Range("B276").Value = 0
$ws.Range("C276").Value = 1
$ws.Range("D276").Value = 16.17337861879347

$ws.Range("A277").PasteSpecial(-4122)
$ws.Range("A277").Value = 44351
$ws.Range("B277").Value = 0
$ws.Range("C277").Value = 1
$ws.Range("D277").Value = 16.17337861879347

$ws.Range("A278").PasteSpecial(-4122)
$ws.Range("A278").Value = 44352
$ws.Range("B278").Value = 0
$ws.Range("C278").Value = 0
$ws.Range("D278").Value = 0

$ws.Range("A279").PasteSpecial(-4122)
$ws.Range("A279").Value = 44353
$ws.Range("B279").Value = 0
$ws.Range("C279").Value = 0
$ws.Range("D279").Value = 0

$ws.Range("A280").PasteSpecial(-4122)
$ws.Range("A280").Value = 44354
$ws.Range("B280").Value = 0
$ws.Range("C280").Value = 0
$ws.Range("D280").Value = 0

$ws.Range("A281").PasteSpecial(-4122)
$ws.Range("A281").Value = 44355
$ws.Range("B281").Value = 1
$ws.Range("C281").Value = 1
$ws.Range("D281").Value = 16.17337861879347

$ws.Range("A282").PasteSpecial(-4122)
$ws.Range("A282").Value = 44356
$ws.Range("B282").Value = 0
$ws.Range("C282").Value = 1
$ws.Range("D282").Value = 16.17337861879347

$ws.Range("A283").PasteSpecial(-4122)
$ws.Range("A283").Value = 44357
$ws.Range("B283").Value = 0
$ws.Range("C283").Value = 1
$ws.Range("D283").Value = 16.17337861879347

$ws.Range("A284").PasteSpecial(-4122)
$ws.Range("A284").Value = 44358
$ws.Range("B284").Value = 0
$ws.Range("C284").Value = 1
$ws.Range("D284").Value = 16.17337861879347

$ws.Range("A285").PasteSpecial(-4122)
$ws.Range("A285").Value = 44359
$ws.Range("B285").Value = 0
$ws.Range("C285").Value = 1
$ws.Range("D285").Value = 16.17337861879347

$ws.Range("A286").PasteSpecial(-4122)
$ws.Range("A286").Value = 44360
$ws.Range("B286").Value = 0
$ws.Range("C286").Value = 1
$ws.Range("D286").Value = 16.17337861879347

$ws.Range("A287").PasteSpecial(-4122)
$ws.Range("A287").Value = 44361
$ws.Range("B287").Value = 0
$ws.Range("C287").Value = 1
$ws.Range("D287").Value = 16.17337861879347

$ws.Range("A288").PasteSpecial(-4122)
$ws.Range("A288").Value = 44362
$ws.Range("B288").Value = 0
$ws.Range("C288").Value = 0
$ws.Range("D288").Value = 0

$ws.Range("A289").PasteSpecial(-4122)
$ws.Range("A289").Value = 44363
$ws.Range("B289").Value = 0
$ws.Range("C289").Value = 0
$ws.Range("D289").Value = 0

$ws.Range("A290").PasteSpecial(-4122)
$ws.Range("A290").Value = 44364
$ws.Range("B290").Value = 2
$ws.Range("C290").Value = 2
$ws.Range("D290").Value = 32.34675723758694

$ws.Range("A291").PasteSpecial(-4122)
$ws.Range("A291").Value = 44365
$ws.Range("B291").Value = 0
$ws.Range("C291").Value = 2
$ws.Range("D291").Value = 32.34675723758694

$ws.Range("A292").PasteSpecial(-4122)
$ws.Range("A292").Value = 44366
$ws.Range("B292").Value = 0
$ws.Range("C292").Value = 2
$ws.Range("D292").Value = 32.34675723758694

$ws.Range("A293").PasteSpecial(-4122)
$ws.Range("A293").Value = 44367
$ws.Range("B293").Value = 0
$ws.Range("C293").Value = 2
$ws.Range("D293").Value = 32.34675723758694

$ws.Range("A294").PasteSpecial(-4122)
$ws.Range("A294").Value = 44368
$ws.Range("B294").Value = 0
$ws.Range("C294").Value = 2
$ws.Range("D294").Value = 32.34675723758694

$ws.Range("A295").PasteSpecial(-4122)
$ws.Range("A295").Value = 44369
$ws.Range("B295").Value = 0
$ws.Range("C295").Value = 2
$ws.Range("D295").Value = 32.34675723758694

$ws.Range("A296").PasteSpecial(-4122)
$ws.Range("A296").Value = 44370
$ws.Range("B296").Value = 0
$ws.Range("C296").Value = 2
$ws.Range("D296").Value = 32.34675723758694

$ws.Range("A297").PasteSpecial(-4122)
$ws.Range("A297").Value = 44371
$ws.Range("B297").Value = 0
$ws.Range("C297").Value = 0
$ws.Range("D297").Value = 0

$ws.Range("A298").PasteSpecial(-4122)
$ws.Range("A298").Value = 44372
$ws.Range("B298").Value = 0
$ws.Range("C298").Value = 0
$ws.Range("D298").Value = 0

$ws.Range("A299").PasteSpecial(-4122)
$ws.Range("A299").Value = 44373
$ws.Range("B299").Value = 0
$ws.Range("C299").Value = 0
$ws.Range("D299").Value = 0

$ws.Range("A300").PasteSpecial(-4122)
$ws.Range("A300").Value = 44374
$ws.Range("B300").Value = 0
$ws.Range("C300").Value = 0
$ws.Range("D300").Value = 0

$ws.Range("A301").PasteSpecial(-4122)
$ws.Range("A301").Value = 44375
$ws.Range("B301").Value = 1
$ws.Range("C301").Value = 1
$ws.Range("D301").Value = 16.17337861879347

$excel.CutCopyMode = 0